$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in "Week 9" data -------------------------------------------------
# Header: give the placeholder "Week 9" header its real date range.
$ws.Range("J1").Value2 = "Week 9 -- June27 - July 3"

# Weekly totals for each athlete (row 10 / Brandon Greife stays 0 - unchanged).
$ws.Range("J2").Value2  = 141.1
$ws.Range("J3").Value2  = 83.8
$ws.Range("J4").Value2  = 79.2
$ws.Range("J5").Value2  = 274.89999999999998
$ws.Range("J6").Value2  = 321.5
$ws.Range("J7").Value2  = 104.9
$ws.Range("J8").Value2  = 65.900000000000006
$ws.Range("J9").Value2  = 88
$ws.Range("J11").Value2 = 173.7

# Column J now holds the longer header text - widen it to fit (mirrors the
# autofit Excel performs once the cell holds the longer date-range string).
$ws.Columns.Item(10).ColumnWidth = 20.5

# --- Stray formatted-but-empty cells picked up along the way --------------
$ws.Range("H13").NumberFormat = "0.0"
$ws.Range("I13").NumberFormat = "0.0"
$ws.Range("I15").NumberFormat = "0.0"
$ws.Range("I19").NumberFormat = "0.0"
$ws.Range("H23").NumberFormat = "0.0"
$ws.Range("I23").NumberFormat = "0.0"
$ws.Range("I33").NumberFormat = "0.0"

# --- View state -------------------------------------------------------------
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("E14").Select()
